$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update the "time_taken" column (F) on the "data" sheet with refreshed timestamps ---
$timestamps = @(
    "2021-10-05 14:35:42.034428",
    "2021-10-05 14:35:42.034436",
    "2021-10-05 14:35:42.034439",
    "2021-10-05 14:35:42.034441",
    "2021-10-05 14:35:42.034444",
    "2021-10-05 14:35:42.034447",
    "2021-10-05 14:35:42.034450",
    "2021-10-05 14:35:42.034452",
    "2021-10-05 14:35:42.034455",
    "2021-10-05 14:35:42.034458",
    "2021-10-05 14:35:42.034460",
    "2021-10-05 14:35:42.034463"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- Add a new "metadata" worksheet placed right after "data" ---
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

# Header row (B1:G1) - values then formatting copied from the "data" sheet's header style
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $meta.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row 2
$meta.Cells.Item(2, 1).Value = 0
$ws.Range("A2").Copy()
$meta.Cells.Item(2, 1).PasteSpecial(-4122)

$meta.Cells.Item(2, 2).Value = "Severe Combined Immunodeficiency (absent T present B cells)"
$meta.Cells.Item(2, 3).Value = 235

# Keep "1.0" as literal text (not converted to the number 1) without leaving a custom style behind:
# enter it as a text-producing formula, then paste-special "values" over itself to flatten
# it back down to a plain string cell (no <f>, no new number format / style).
$meta.Cells.Item(2, 4).Formula = '="1.0"'
$meta.Cells.Item(2, 4).Copy()
$meta.Cells.Item(2, 4).PasteSpecial(-4163)
$excel.CutCopyMode = $false

$meta.Cells.Item(2, 5).Value = "2021-09-30T07:44:36.985556Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:35:42.030769"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/235/?format=json"

# Keep "data" as the active/selected sheet, matching the original workbook view
$ws.Activate()
